$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (diff: name="Sheet" -> name="Sheet1")
$ws.Name = "Sheet1"

# Clear out the old stray cell data (B10 = 60.7) so the used range resets
$ws.Cells.Clear()

# Write the header row
$headers = @("Date", "Open", "High", "Low", "Close", "Volume", "Dividends", "Stock Splits")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4160    # xlTop
    $cell.Borders.LineStyle = 1        # xlContinuous
    $cell.Borders.Weight = 2           # xlThin
}
